$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append four new data rows (rows 6-9) to the Fleet Space Report sheet,
# extending the used range from A1:M5 to A1:N9. Column N has no header,
# matching the source diff (row 1 is left untouched).

$data = @(
    @("2025-02-07 13:41", "pstg-fa-03", 1,                 102402176,     0,           0, 1,                 107961344,    1,                 5559168,     0,            1111537536204800, 107961344,    0),
    @("2025-02-07 13:41", "pstg-fa-02", 4.13294045762392,  142350152779,  58263161536, 0, 0.999815000270004, 327266564890, 22251.60979597833, 126653250575, 1029782972928, 5566402572320768, 327266564890, 0),
    @("2025-02-07 15:22", "pstg-fa-03", 1,                 102402176,     0,           0, 1,                 107961344,    1,                 5559168,     0,            1111537536204800, 107961344,    0),
    @("2025-02-07 15:22", "pstg-fa-02", 4.134110249036264, 142260835620,  58262194555, 0, 0.9998149961214116, 327178912324, 22257.41077400217, 126655882149, 1029806065664, 5566402572320768, 327178912324, 0)
)

$startRow = 6
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
}
